$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(2, 4).Value = '65.056.98'
$ws.Cells.Item(2, 5).Value = '  +2.00%  '

$ws.Cells.Item(3, 4).Value = '3.393.87'
$ws.Cells.Item(3, 5).Value = '  +1.64%  '

$ws.Cells.Item(4, 5).Value = '  +0.18%  '

$c = $ws.Cells.Item(5, 4)
$c.NumberFormat = "@"
$c.Value = '558.80'
$ws.Cells.Item(5, 5).Value = '  +2.24%  '

$c = $ws.Cells.Item(6, 4)
$c.NumberFormat = "@"
$c.Value = '175.12'
$ws.Cells.Item(6, 5).Value = '  +1.75%  '

$ws.Cells.Item(7, 5).Value = '  +2.12%  '

$ws.Cells.Item(8, 4).Value = '3.383.97'
$ws.Cells.Item(8, 5).Value = '  +1.77%  '

$ws.Cells.Item(9, 5).Value = '  +0.05%  '

$c = $ws.Cells.Item(10, 4)
$c.NumberFormat = "@"
$c.Value = '0.169'
$ws.Cells.Item(10, 5).Value = '  +11.60%  '

$c = $ws.Cells.Item(11, 4)
$c.NumberFormat = "@"
$c.Value = '0.630'
$ws.Cells.Item(11, 5).Value = '  +3.28%  '

$c = $ws.Cells.Item(12, 4)
$c.NumberFormat = "@"
$c.Value = '54.76'
$ws.Cells.Item(12, 5).Value = '  +2.25%  '

$c = $ws.Cells.Item(13, 4)
$c.NumberFormat = "@"
$c.Value = '0.0000279'
$ws.Cells.Item(13, 5).Value = '  +5.37%  '

$c = $ws.Cells.Item(14, 4)
$c.NumberFormat = "@"
$c.Value = '9.13'
$ws.Cells.Item(14, 5).Value = '  +3.01%  '

$ws.Cells.Item(15, 4).Value = '3.940.34'
$ws.Cells.Item(15, 5).Value = '  +8.37%  '

$c = $ws.Cells.Item(16, 4)
$c.NumberFormat = "@"
$c.Value = '18.38'
$ws.Cells.Item(16, 5).Value = '  +1.88%  '

$ws.Cells.Item(17, 2).Value = 'TRON'
$ws.Cells.Item(17, 3).Value = 'https://coinranking.com/coin/qUhEFk1I61atv+tron-trx'
$c = $ws.Cells.Item(17, 4)
$c.NumberFormat = "@"
$c.Value = '0.119'
$ws.Cells.Item(17, 5).Value = '  +1.86%  '

$ws.Cells.Item(18, 2).Value = 'WrappedEther'
$ws.Cells.Item(18, 3).Value = 'https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth'
$ws.Cells.Item(18, 4).Value = '3.393.48'
$ws.Cells.Item(18, 5).Value = '  +3.03%  '

$ws.Cells.Item(19, 4).Value = '65.043.39'
$ws.Cells.Item(19, 5).Value = '  +2.23%  '

$c = $ws.Cells.Item(20, 4)
$c.NumberFormat = "@"
$c.Value = '11.85'
$ws.Cells.Item(20, 5).Value = '  +1.42%  '

$c = $ws.Cells.Item(21, 4)
$c.NumberFormat = "@"
$c.Value = '0.993'
$ws.Cells.Item(21, 5).Value = '  +1.88%  '

$c = $ws.Cells.Item(22, 4)
$c.NumberFormat = "@"
$c.Value = '473.57'
$ws.Cells.Item(22, 5).Value = '  +15.62%  '

$ws.Cells.Item(23, 5).Value = '  +13.85%  '

$c = $ws.Cells.Item(24, 4)
$c.NumberFormat = "@"
$c.Value = '4.13'
$ws.Cells.Item(24, 5).Value = '  +2.76%  '

$c = $ws.Cells.Item(25, 4)
$c.NumberFormat = "@"
$c.Value = '87.30'
$ws.Cells.Item(25, 5).Value = '  +5.24%  '

$c = $ws.Cells.Item(26, 4)
$c.NumberFormat = "@"
$c.Value = '13.47'
$ws.Cells.Item(26, 5).Value = '  -1.60%  '

$c = $ws.Cells.Item(27, 4)
$c.NumberFormat = "@"
$c.Value = '2.90'
$ws.Cells.Item(27, 5).Value = '  +6.47%  '

$c = $ws.Cells.Item(28, 4)
$c.NumberFormat = "@"
$c.Value = '10.87'
$ws.Cells.Item(28, 5).Value = '  +3.16%  '

$c = $ws.Cells.Item(29, 4)
$c.NumberFormat = "@"
$c.Value = '8.79'
$ws.Cells.Item(29, 5).Value = '  +2.39%  '

$c = $ws.Cells.Item(30, 4)
$c.NumberFormat = "@"
$c.Value = '31.19'
$ws.Cells.Item(30, 5).Value = '  +7.44%  '

$c = $ws.Cells.Item(31, 4)
$c.NumberFormat = "@"
$c.Value = '6.69'
$ws.Cells.Item(31, 5).Value = '  +5.18%  '

$c = $ws.Cells.Item(32, 4)
$c.NumberFormat = "@"
$c.Value = '11.54'
$ws.Cells.Item(32, 5).Value = '  +1.75%  '

$c = $ws.Cells.Item(33, 4)
$c.NumberFormat = "@"
$c.Value = '61.96'
$ws.Cells.Item(33, 5).Value = '  +6.99%  '

$c = $ws.Cells.Item(34, 4)
$c.NumberFormat = "@"
$c.Value = '572.73'
$ws.Cells.Item(34, 5).Value = '  -0.56%  '

$ws.Cells.Item(35, 5).Value = '  +2.05%  '

$c = $ws.Cells.Item(36, 4)
$c.NumberFormat = "@"
$c.Value = '1.00'

$ws.Cells.Item(37, 5).Value = '  -4.81%  '

$ws.Cells.Item(38, 5).Value = '  +3.83%  '

$ws.Cells.Item(39, 2).Value = 'InjectiveProtocol'
$ws.Cells.Item(39, 3).Value = 'https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj'
$c = $ws.Cells.Item(39, 4)
$c.NumberFormat = "@"
$c.Value = '35.77'
$ws.Cells.Item(39, 5).Value = '  +1.75%  '

$ws.Cells.Item(40, 2).Value = 'PEPE'
$ws.Cells.Item(40, 3).Value = 'https://coinranking.com/coin/03WI8NQPF+pepe-pepe'
$ws.Cells.Item(40, 4).Value = '0.0₃0759'
$ws.Cells.Item(40, 5).Value = '  +2.89%  '

$c = $ws.Cells.Item(41, 4)
$c.NumberFormat = "@"
$c.Value = '0.371'
$ws.Cells.Item(41, 5).Value = '  +1.40%  '

$ws.Cells.Item(42, 4).Value = '3.091.40'
$ws.Cells.Item(42, 5).Value = '  -1.16%  '

$ws.Cells.Item(43, 5).Value = '  +0.22%  '

$c = $ws.Cells.Item(44, 4)
$c.NumberFormat = "@"
$c.Value = '2.85'
$ws.Cells.Item(44, 5).Value = '  +2.02%  '

$ws.Cells.Item(45, 5).Value = '  +4.20%  '

$ws.Cells.Item(46, 5).Value = '  +6.10%  '

$c = $ws.Cells.Item(47, 4)
$c.NumberFormat = "@"
$c.Value = '2.47'
$ws.Cells.Item(47, 5).Value = '  +2.44%  '

$c = $ws.Cells.Item(48, 4)
$c.NumberFormat = "@"
$c.Value = '3.17'
$ws.Cells.Item(48, 5).Value = '  -2.21%  '

$ws.Cells.Item(49, 5).Value = '  +0.60%  '

$c = $ws.Cells.Item(50, 4)
$c.NumberFormat = "@"
$c.Value = '137.65'
$ws.Cells.Item(50, 5).Value = '  +4.19%  '

$c = $ws.Cells.Item(51, 4)
$c.NumberFormat = "@"
$c.Value = '8.31'
$ws.Cells.Item(51, 5).Value = '  +3.36%  '
